$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "all": append 2020-05-16 (serial 43967) row at row 39.
# Rows.Insert() pushes the old row 39 (footer note) down to row 40 and
# copies formatting from the row above, matching the existing pattern.
# ------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("all")
$wsAll.Activate()
$wsAll.Rows.Item(39).Insert()
$wsAll.Range("A39").Value = 43967
$wsAll.Range("B39").Value = 282
$wsAll.Range("C39").Value = 278
$wsAll.Range("D39").Value = 63
$wsAll.Range("E39").Value = 53
$wsAll.Range("F39").Value = 10
$wsAll.Range("G39").Value = 11
$wsAll.Range("H39").Value = 204
$wsAll.Range("A39").Select()

# ------------------------------------------------------------------
# Sheet "kobe": the latest existing row (93) gets revised totals, then
# a brand-new row (94) is inserted for 2020-05-16 (serial 43967).
# ------------------------------------------------------------------
$wsKobe = $wb.Worksheets.Item("kobe")
$wsKobe.Activate()
$wsKobe.Range("D93").Value = 1
$wsKobe.Range("E93").Value = 282
$wsKobe.Range("F93").Value = 58

$wsKobe.Rows.Item(94).Insert()
$wsKobe.Range("A94").Value = 43967
$wsKobe.Range("B94").Value = 0
$wsKobe.Range("C94").Value = 2813
$wsKobe.Range("D94").Value = 0
$wsKobe.Range("E94").Value = 282
$wsKobe.Range("F94").Value = 58
$wsKobe.Range("G94").Value = 49
$wsKobe.Range("H94").Value = 9
$wsKobe.Range("I94").Value = 11
$wsKobe.Range("J94").Value = 195
$wsKobe.Range("A94").Select()

# ------------------------------------------------------------------
# Sheet "other": append 2020-05-16 (serial 43967) row at row 69, values
# unchanged from the previous day (no new non-Kobe cases).
# ------------------------------------------------------------------
$wsOther = $wb.Worksheets.Item("other")
$wsOther.Activate()
$wsOther.Rows.Item(69).Insert()
$wsOther.Range("A69").Value = 43967
$wsOther.Range("B69").Value = 0
$wsOther.Range("C69").Value = 14
$wsOther.Range("D69").Value = 5
$wsOther.Range("E69").Value = 4
$wsOther.Range("F69").Value = 1
$wsOther.Range("G69").Value = 0
$wsOther.Range("H69").Value = 9
$wsOther.Range("A68").Select()

# Restore "all" as the active sheet/tab, matching the original workbook.
$wsAll.Activate()
$wsAll.Range("A39").Select()
